$d = $word.ActiveDocument

function Remove-ParagraphRange($doc, $firstDeleteIndex, $firstKeepIndex) {
    # Deletes paragraphs [firstDeleteIndex .. firstKeepIndex) using ORIGINAL
    # (pre-deletion) 1-based paragraph indices of $doc at call time.
    $s = $doc.Paragraphs($firstDeleteIndex).Range.Start
    $e = $doc.Paragraphs($firstKeepIndex).Range.Start
    $doc.Range($s, $e).Delete()
}

# ---------------------------------------------------------------------------
# 1) Remove the bolded section headers (and the whole paragraphs that are no
#    longer needed) from the bottom of the document upward, so that higher
#    paragraph indices referenced below stay valid while we work.
# ---------------------------------------------------------------------------

# "Thank you for your continued support..." paragraph + its trailing blank
Remove-ParagraphRange $d 30 32

# "Gratitude" bold header + its trailing blank
Remove-ParagraphRange $d 26 28

# "Looking Forward" bold header + its trailing blank
Remove-ParagraphRange $d 22 24

# "Challenges and Learning" bold header + its trailing blank
# (keep the paragraph that follows - its text will be replaced below)
Remove-ParagraphRange $d 18 20

# "Achievements and Progress" bold header, its trailing blank, the
# "African Private Sector Forum..." paragraph, and the blank after it
Remove-ParagraphRange $d 10 14

# ---------------------------------------------------------------------------
# 2) Replace the remaining paragraph texts with the new copy. After the
#    deletions above the document has settled into its final 24-paragraph
#    shape, so we can address paragraphs by their new indices.
# ---------------------------------------------------------------------------

$d.Paragraphs(4).Range.Text = "Dear Stakeholders and Partners,"

$d.Paragraphs(6).Range.Text = "As we reflect on the year 2023, the IKEA Foundation has continued to strive toward our mission of creating a better everyday life for the many people, particularly those living in vulnerable conditions due to poverty and climate change. This year has been marked by significant achievements and some challenges, which have tested our resilience and strengthened our resolve."

$d.Paragraphs(8).Range.Text = 'We are proud to report substantial progress across our various projects. Notably, our initiative "Donation: Africa Private Sector Forum on Forced Displacement" has made strides in promoting Decent Work and Economic Growth (Goal 8 of the United Nations)[x]. This project, based in Nairobi, Kenya, aims to enhance employment opportunities within refugee communities, a critical step towards economic inclusion and stability.'

$d.Paragraphs(10).Range.Text = 'Our commitment to combating climate change has been exemplified by the "Funders Table" project, which aims to reduce carbon emissions by 30% within the funded projects[x]. This global initiative underscores our dedication to United Nations'' Goal 13, taking urgent action to combat climate change and its impacts.'

$d.Paragraphs(12).Range.Text = 'The "Just Transition Fund" is another cornerstone project that has been instrumental in reducing carbon emissions and increasing renewable energy use across multiple countries, including South Africa, Vietnam, and Indonesia[x]. These efforts align with our strategic focus on bolstering climate resilience and supporting sustainable energy transitions.'

$d.Paragraphs(14).Range.Text = 'In addressing the urgent needs brought about by natural disasters, our response to the Türkiye and northern Syria earthquake through "Donation: Türkiye and northern Syria earthquake response" provided critical medical aid and reconstruction support, reflecting our commitment to Sustainable Cities and Communities (UN Goal 11)[x].'

$d.Paragraphs(16).Range.Text = "As we look to the future, our vision is clear. We aim to deepen our impact, expand our reach, and continue to innovate in our approach to philanthropy. The challenges of poverty and climate change are intertwined and complex, but with the continued support and collaboration from our partners and stakeholders, we are optimistic about what we can achieve together."

$d.Paragraphs(18).Range.Text = "We express our deepest gratitude to all our supporters who share our vision and commitment. Your unwavering support fuels our efforts to make a significant and lasting impact."

$d.Paragraphs(20).Range.Text = "Sincerely,"

$d.Paragraphs(23).Range.Text = "Chairperson/President, IKEA Foundation"

$d.Paragraphs(24).Range.Text = "[Date] 2024"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
